# Applies the ExecutorFrameWork.docx edit:
#   1. "Fork/Joi" + "n" (two runs) -> "Fork/Join" (one run).
#   2. The "_GoBack" bookmark moves from the end of the document to the
#      middle of the "...get method will throw an..." sentence, splitting
#      the run that contains it into "...get m" | bookmark | "ethod...".

$d = $word.ActiveDocument

# --- 1. Merge "Fork/Joi" + "n" into a single "Fork/Join" run -------------
# The two runs share identical run formatting, so replacing the combined
# "Fork/Join" text with itself collapses them into a single run.
$d.Content.Find.Execute("Fork/Join", $true, $false, $false, $false, $false, $true, 1, $false, "Fork/Join", 2)

# --- 2. Relocate the "_GoBack" bookmark into the exception sentence ------
$prefix = "If there is an exception when executing the task, the call to get m"
$sentence = $prefix + "ethod will throw an "

$rng = $d.Content
$rng.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($rng.Find.Found) {
    $splitPoint = $rng.Start + $prefix.Length
    $point = $d.Range($splitPoint, $splitPoint)

    # Bookmarks.Add with a name that already exists relocates that
    # bookmark: it disappears from its old spot (the empty paragraph right
    # before the final sectPr) and reappears here, splitting the run in
    # two exactly like the target diff.
    $d.Bookmarks.Add("_GoBack", $point)
}
